$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '51.757.76'
Set-TextValue $ws.Range("E2") '  +0.74%  '

Set-TextValue $ws.Range("D3") '2.825.60'
Set-TextValue $ws.Range("E3") '  +1.87%  '

Set-TextValue $ws.Range("E4") '  +0.04%  '

Set-TextValue $ws.Range("D5") '351.24'
Set-TextValue $ws.Range("E5") '  -0.77%  '

Set-TextValue $ws.Range("D6") '113.06'
Set-TextValue $ws.Range("E6") '  +4.40%  '

Set-TextValue $ws.Range("E7") '  +1.22%  '

Set-TextValue $ws.Range("E8") '  +0.08%  '

Set-TextValue $ws.Range("E9") '  +6.04%  '

Set-TextValue $ws.Range("D10") '40.11'
Set-TextValue $ws.Range("E10") '  +1.16%  '

Set-TextValue $ws.Range("E11") '  -0.86%  '

Set-TextValue $ws.Range("D12") '0.0851'
Set-TextValue $ws.Range("E12") '  +2.21%  '

Set-TextValue $ws.Range("D13") '19.99'
Set-TextValue $ws.Range("E13") '  -0.26%  '

Set-TextValue $ws.Range("D14") '7.79'
Set-TextValue $ws.Range("E14") '  +3.45%  '

Set-TextValue $ws.Range("D15") '3.272.75'
Set-TextValue $ws.Range("E15") '  +2.07%  '

Set-TextValue $ws.Range("E16") '  +5.13%  '

Set-TextValue $ws.Range("D17") '2.813.37'
Set-TextValue $ws.Range("E17") '  +1.64%  '

Set-TextValue $ws.Range("D18") '51.797.07'
Set-TextValue $ws.Range("E18") '  +0.97%  '

Set-TextValue $ws.Range("E19") '  +11.43%  '

Set-TextValue $ws.Range("D20") '7.60'
Set-TextValue $ws.Range("E20") '  -0.30%  '

Set-TextValue $ws.Range("D21") '13.36'
Set-TextValue $ws.Range("E21") '  +1.53%  '

Set-TextValue $ws.Range("E22") '  +1.36%  '

Set-TextValue $ws.Range("D23") '70.59'
Set-TextValue $ws.Range("E23") '  +1.24%  '

Set-TextValue $ws.Range("D24") '269.06'
Set-TextValue $ws.Range("E24") '  +1.20%  '

Set-TextValue $ws.Range("D25") '2.77'
Set-TextValue $ws.Range("E25") '  +2.43%  '

Set-TextValue $ws.Range("D26") '26.27'
Set-TextValue $ws.Range("E26") '  +1.08%  '

Set-TextValue $ws.Range("E27") '  -0.01%  '

Set-TextValue $ws.Range("E28") '  -0.14%  '

Set-TextValue $ws.Range("D29") '38.97'
Set-TextValue $ws.Range("E29") '  +7.42%  '

Set-TextValue $ws.Range("E30") '  +3.67%  '

Set-TextValue $ws.Range("D31") '2.25'
Set-TextValue $ws.Range("E31") '  +2.09%  '

Set-TextValue $ws.Range("D32") '52.75'
Set-TextValue $ws.Range("E32") '  +1.86%  '

Set-TextValue $ws.Range("E33") '  +1.12%  '

Set-TextValue $ws.Range("D34") '0.0458'
Set-TextValue $ws.Range("E34") '  +3.46%  '

Set-TextValue $ws.Range("E35") '  +8.80%  '

Set-TextValue $ws.Range("D36") '5.65'
Set-TextValue $ws.Range("E36") '  +2.53%  '

Set-TextValue $ws.Range("D37") '0.999'
Set-TextValue $ws.Range("E37") '  +0.01%  '

Set-TextValue $ws.Range("D38") '19.07'
Set-TextValue $ws.Range("E38") '  +4.87%  '

Set-TextValue $ws.Range("E39") '  +2.67%  '

Set-TextValue $ws.Range("E40") '  +2.64%  '

Set-TextValue $ws.Range("E41") '  +1.49%  '

Set-TextValue $ws.Range("E42") '  +1.09%  '

Set-TextValue $ws.Range("D43") '121.76'
Set-TextValue $ws.Range("E43") '  +0.93%  '

Set-TextValue $ws.Range("E44") '  +1.27%  '

Set-TextValue $ws.Range("D45") '22.10'
Set-TextValue $ws.Range("E45") '  +0.71%  '

Set-TextValue $ws.Range("D46") '3.50'
Set-TextValue $ws.Range("E46") '  +7.98%  '

Set-TextValue $ws.Range("B47") 'Maker'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D47") '2.174.32'
Set-TextValue $ws.Range("E47") '  +3.68%  '

Set-TextValue $ws.Range("B48") 'ApeXProtocol'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D48") '2.49'
Set-TextValue $ws.Range("E48") '  +7.52%  '

Set-TextValue $ws.Range("B49") 'TheGraph'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D49") '0.243'
Set-TextValue $ws.Range("E49") '  +26.22%  '

Set-TextValue $ws.Range("B50") 'SEI'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
Set-TextValue $ws.Range("D50") '0.996'
Set-TextValue $ws.Range("E50") '  +10.39%  '

Set-TextValue $ws.Range("D51") '5.53'
Set-TextValue $ws.Range("E51") '  +2.26%  '
